# Update Efna4/Epha3 LR-pair table with refreshed TPM-based values.
# Sending/target cluster set changes from {ECs,FAPs,MuSCs} x {ECs,FAPs,MuSCs}+Resolving-Mac(sender only)
# to {ECs,FAPs,MuSCs} x {ECs,FAPs,MuSCs,Resolving-Mac} (Resolving-Mac becomes target-only).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = [object[]]@("ECs", "Efna4", "Epha3", "ECs", 3, 1, 1.008000333333333, 3.024001, 0.3525296793986107, 0.3525296793986107, 2, 0.6666666666666666, 0.07579599999999999, 0.227388, 0.001780200955210419, 0.001780200955210419, 0.07640239326533332, 0.687621539388, 0.0006275736720054295, 0.0006275736720054294)
$row3 = [object[]]@("ECs", "Efna4", "Epha3", "FAPs", 3, 1, 1.008000333333333, 3.024001, 0.3525296793986107, 0.3525296793986107, 3, 1, 40.12734033333333, 120.382021, 0.9424604146848589, 0.9424604146848587, 40.44837243178011, 364.035351886021, 0.332245267834735, 0.332245267834735)
$row4 = [object[]]@("ECs", "Efna4", "Epha3", "MuSCs", 3, 1, 1.008000333333333, 3.024001, 0.3525296793986107, 0.3525296793986107, 3, 1, 2.344072666666666, 7.032217999999999, 0.05505462557763778, 0.05505462557763778, 2.362826029357556, 21.265434264218, 0.0194083895042952, 0.0194083895042952)
$row5 = [object[]]@("ECs", "Efna4", "Epha3", "Resolving-Mac", 3, 1, 1.008000333333333, 3.024001, 0.3525296793986107, 0.3525296793986107, 2, 0.6666666666666666, 0.03000666666666667, 0.09002, 0.0007047587822930054, 0.0007047587822930053, 0.03024673000222223, 0.27222057002, 0.0002484483875751085, 0.0002484483875751084)
$row6 = [object[]]@("FAPs", "Efna4", "Epha3", "ECs", 3, 1, 1.537538333333333, 4.612615, 0.5377259091975243, 0.5377259091975243, 2, 0.6666666666666666, 0.07579599999999999, 0.227388, 0.001780200955210419, 0.001780200955210419, 0.1165392555133333, 1.04885329962, 0.0009572601771948235, 0.0009572601771948235)
$row7 = [object[]]@("FAPs", "Efna4", "Epha3", "FAPs", 3, 1, 1.537538333333333, 4.612615, 0.5377259091975243, 0.5377259091975243, 3, 1, 40.12734033333333, 120.382021, 0.9424604146848589, 0.9424604146848587, 61.69732397721277, 555.2759157949149, 0.5067853833690915, 0.5067853833690914)
$row8 = [object[]]@("FAPs", "Efna4", "Epha3", "MuSCs", 3, 1, 1.537538333333333, 4.612615, 0.5377259091975243, 0.5377259091975243, 3, 1, 2.344072666666666, 7.032217999999999, 0.05505462557763778, 0.05505462557763778, 3.604101581118889, 32.43691423006999, 0.02960429859426455, 0.02960429859426455)
$row9 = [object[]]@("FAPs", "Efna4", "Epha3", "Resolving-Mac", 3, 1, 1.537538333333333, 4.612615, 0.5377259091975243, 0.5377259091975243, 2, 0.6666666666666666, 0.03000666666666667, 0.09002, 0.0007047587822930054, 0.0007047587822930053, 0.04613640025555556, 0.4152276023, 0.0003789670569734464, 0.0003789670569734464)
$row10 = [object[]]@("MuSCs", "Efna4", "Epha3", "ECs", 2, 0.6666666666666666, 0.313796, 0.941388, 0.1097444114038651, 0.1097444114038651, 2, 0.6666666666666666, 0.07579599999999999, 0.227388, 0.001780200955210419, 0.001780200955210419, 0.023784481616, 0.214060334544, 0.0001953671060101657, 0.0001953671060101657)
$row11 = [object[]]@("MuSCs", "Efna4", "Epha3", "FAPs", 2, 0.6666666666666666, 0.313796, 0.941388, 0.1097444114038651, 0.1097444114038651, 3, 1, 40.12734033333333, 120.382021, 0.9424604146848589, 0.9424604146848587, 12.59179888723867, 113.326189985148, 0.1034297634810324, 0.1034297634810324)
$row12 = [object[]]@("MuSCs", "Efna4", "Epha3", "MuSCs", 2, 0.6666666666666666, 0.313796, 0.941388, 0.1097444114038651, 0.1097444114038651, 3, 1, 2.344072666666666, 7.032217999999999, 0.05505462557763778, 0.05505462557763778, 0.7355606265093333, 6.620045638583999, 0.006041937479078033, 0.006041937479078032)
$row13 = [object[]]@("MuSCs", "Efna4", "Epha3", "Resolving-Mac", 2, 0.6666666666666666, 0.313796, 0.941388, 0.1097444114038651, 0.1097444114038651, 2, 0.6666666666666666, 0.03000666666666667, 0.09002, 0.0007047587822930054, 0.0007047587822930053, 0.009415971973333334, 0.08474374776, 0.000077343337744450565605057457, 0.000077343337744450538500003145)

$allRows = [object[]]@($row2, $row3, $row4, $row5, $row6, $row7, $row8, $row9, $row10, $row11, $row12, $row13)

$rowNum = 2
foreach ($row in $allRows) {
    $colNum = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowNum, $colNum).Value = $val
        $colNum = $colNum + 1
    }
    $rowNum = $rowNum + 1
}
